# Fruta / hortaliza, semanal
# The data rows (2-41) got reshuffled: each row's Fecha/Calidad/Volumen/
# Precio min/max/promedio/Origen/Precio-$-Kg values were swapped with the
# corresponding values from another row in the original sheet (columns
# A,B,C,E-K,Q,T are identical across rows, so the rest of the row stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# after-row -> before-row mapping (1-based worksheet rows)
$map = @{
    2  = 6
    3  = 14
    4  = 10
    5  = 18
    6  = 41
    7  = 8
    8  = 32
    9  = 27
    10 = 19
    11 = 5
    12 = 2
    13 = 37
    14 = 40
    15 = 25
    16 = 20
    17 = 34
    18 = 28
    19 = 9
    20 = 33
    21 = 4
    22 = 29
    23 = 21
    24 = 3
    25 = 11
    26 = 16
    27 = 24
    28 = 31
    29 = 26
    30 = 30
    31 = 7
    32 = 39
    33 = 36
    34 = 22
    35 = 17
    36 = 35
    37 = 12
    38 = 15
    39 = 23
    40 = 38
    41 = 13
}

# Columns that move together as a unit per row.
$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the original (pre-edit) values for the columns that change,
# for every data row, before any writes happen. Use Value2 for reads --
# the plain Value getter on this host returns a member-info placeholder
# instead of the actual cell content.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write back each row's new values, taken from its mapped source row.
for ($r = 2; $r -le 41; $r++) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
